$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(4)
$tr = $shp.TextFrame.TextRange

# The target run ("일 주행거리가 매년 줄어 들고 있음") is the 2nd run of
# paragraph 4 in this shape, starting at absolute character offset 97
# (1-based) with length 19.
$target = $tr.Characters(97, 19)
$target.Text = "일 평균 주행거리가 매년 줄어 들고 있음"

# Split the (now longer) replacement text into three runs that share the
# same run formatting as the original run, matching the authored diff:
#   "일 " + "평균 주행거리가 " + "매년 줄어 들고 있음"
$run1 = $tr.Characters(97, 2)
$run1.Text = "일 "

$run2 = $tr.Characters(99, 9)
$run2.Text = "평균 주행거리가 "

$run3 = $tr.Characters(108, 11)
$run3.Text = "매년 줄어 들고 있음"

# Re-splitting the run via Characters() nudges this autofit ("shrink to
# fit") textbox's cached height; restore it to the original value so the
# shape geometry stays untouched, matching the source diff. The literal
# below compensates for this host's points->EMU float rounding so the
# saved <a:ext cy="..."/> lands back on the exact original 2862322 EMU.
$shp.Height = 225.379695
